# Weekly update: add this week's two new Alcachofa price records (Vega Modelo
# de Temuco) at the top of the data block, pushing the existing history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the current first data-history row (88),
# shifting every existing record (old rows 88-102) down to rows 90-104.
$ws.Rows.Item(88).Insert()
$ws.Rows.Item(89).Insert()

# New row 88: Alcachofa "Española" entry for 2021-09-10 (serial 44449)
$ws.Range("A88").Value = 10
$ws.Range("B88").Value = "Vega Modelo de Temuco"
$ws.Range("C88").Value = "La Araucanía"
$ws.Range("D88").Value = 44449
$ws.Range("E88").Value = 9
$ws.Range("F88").Value = 100112013
$ws.Range("G88").Value = "Alcachofa"
$ws.Range("H88").Value = "Española"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 175
$ws.Range("K88").Value = 14000
$ws.Range("L88").Value = 15000
$ws.Range("M88").Value = 14543
$ws.Range("N88").Value = "$/caja 30 unidades"
$ws.Range("O88").Value = "Provincia de Limarí"
$ws.Range("P88").Value = 485
$ws.Range("Q88").Value = 30
$ws.Range("R88").Value = "Hortaliza"

# New row 89: Alcachofa "Madrigal" entry for 2021-09-10 (serial 44449)
$ws.Range("A89").Value = 10
$ws.Range("B89").Value = "Vega Modelo de Temuco"
$ws.Range("C89").Value = "La Araucanía"
$ws.Range("D89").Value = 44449
$ws.Range("E89").Value = 9
$ws.Range("F89").Value = 100112013
$ws.Range("G89").Value = "Alcachofa"
$ws.Range("H89").Value = "Madrigal"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 65
$ws.Range("K89").Value = 14000
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = 14000
$ws.Range("N89").Value = "$/caja 40 unidades"
$ws.Range("O89").Value = "Provincia de Limarí"
$ws.Range("P89").Value = 350
$ws.Range("Q89").Value = 40
$ws.Range("R89").Value = "Hortaliza"
